$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "FilesTab" row (row 4): the "query" column (B4) gets a rewritten Cypher query ---
$newFilesQuery = 'MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE f.file_type in [''TXT'']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '''') as `File Name`,
    coalesce(s.study_name,'''') as `Study Name`,
    coalesce(s.phs_accession,'''') as `Accession`,
    coalesce(p.participant_id, '''') as `Participant ID`,
    coalesce(samp.sample_id, '''') as `Sample ID`,
    coalesce(f.file_type, '''') as `File Type`
ORDER BY f.file_name limit 100'
$ws.Cells.Item(4, 2).Value = $newFilesQuery

# --- "StatQuery" column (C) for all three data rows is replaced with a new combined query ---
$newStatQuery = 'CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE  f.file_type in [''TXT'']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE f.file_type in [''TXT'']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE  f.file_type in [''TXT'']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`'
$ws.Cells.Item(2, 3).Value = $newStatQuery
$ws.Cells.Item(3, 3).Value = $newStatQuery
$ws.Cells.Item(4, 3).Value = $newStatQuery

# The longer query text needs more vertical room; rows grow to Excel's maximum row height
$ws.Rows("2:2").RowHeight = 409.5
$ws.Rows("3:3").RowHeight = 409.5
$ws.Rows("4:4").RowHeight = 409.5

# Move the active selection to C5
$ws.Range("C5").Select() | Out-Null
